# Update the "Experimental results" tables on slide 13 (Numerical/Oversample
# experiment) and slide 14 (the companion experiment) to:
#   - un-bold the previously-highlighted "All Features / with standardization"
#     score for the Logistic Regression row (it is no longer the best score),
#   - fill in the previously-missing "All Features / with standardization"
#     score for the Decision Tree row and bold it (it is now the best score),
#   - center the row-label cells vertically in the second (AUC ROC) table on
#     slide 13.

$p = $ppt.ActivePresentation

# ---- Slide 13 ----
$slide13 = $p.Slides.Item(13)

# Table 1 ("Accuracy"): shape 3
$accuracyTable = $slide13.Shapes.Item(3).Table

$accLrPara = $accuracyTable.Cell(3, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$accLrPara.Font.Bold = $false

$accDtPara = $accuracyTable.Cell(4, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$accDtPara.Text = "0.9789"
$accDtPara.Font.Bold = $true

# Table 2 ("AUC ROC"): shape 5
$aucTable = $slide13.Shapes.Item(5).Table

$aucLrPara = $aucTable.Cell(3, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$aucLrPara.Font.Bold = $false

$aucDtPara = $aucTable.Cell(4, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$aucDtPara.Text = "0.845"
$aucDtPara.Font.Bold = $true

# Vertically center the row-label cells (column 1) for the 3 model rows.
$aucTable.Cell(3, 1).Shape.TextFrame.VerticalAnchor = 3
$aucTable.Cell(4, 1).Shape.TextFrame.VerticalAnchor = 3
$aucTable.Cell(5, 1).Shape.TextFrame.VerticalAnchor = 3

# ---- Slide 14 ----
$slide14 = $p.Slides.Item(14)

# Table 1 ("F1-Score"): shape 3
$f1Table = $slide14.Shapes.Item(3).Table

$f1LrPara = $f1Table.Cell(3, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$f1LrPara.Font.Bold = $false

$f1DtPara = $f1Table.Cell(4, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$f1DtPara.Text = "0.7438"
$f1DtPara.Font.Bold = $true

# Table 2 ("P/R"): shape 5
$prTable = $slide14.Shapes.Item(5).Table

$prLrPara = $prTable.Cell(3, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$prLrPara.Font.Bold = $false

$prDtPara = $prTable.Cell(4, 4).Shape.TextFrame.TextRange.Paragraphs(1, 1)
$prDtPara.Text = "0.937/0.616"
$prDtPara.Font.Bold = $true
